# "SP and Comodity path Changes in Customer page"
# Update the short-code values for "Standard Cust" (B9) and "Ad copy" (B12)
# on the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = "HNP9Y"
$ws.Range("B12").Value = "R8"
